$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "São Paulo"
$ws.Range("B8").Value = "26/08/2025 21:04"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "16"
$ws.Range("C8").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "90"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "Nublado"
